$d = $word.ActiveDocument

$replacements = @(
    @("64×23=1472", "64×47=3008"),
    @("39×61=2379", "71×84=5964"),
    @("87×68=5916", "83×58=4814"),
    @("90×36=3240", "65×13=845"),
    @("29×23=667",  "94×18=1692"),
    @("91×74=6734", "29×16=464"),
    @("45×42=1890", "86×24=2064"),
    @("86×32=2752", "95×78=7410"),
    @("36×48=1728", "18×79=1422"),
    @("60×57=3420", "22×80=1760"),
    @("40×98=3920", "69×23=1587"),
    @("67×91=6097", "57×34=1938"),
    @("54×28=1512", "55×59=3245"),
    @("81×21=1701", "38×67=2546"),
    @("19×99=1881", "45×83=3735"),
    @("85×18=1530", "51×21=1071"),
    @("97×61=5917", "14×81=1134"),
    @("52×86=4472", "87×77=6699"),
    @("91×86=7826", "93×76=7068"),
    @("57×78=4446", "89×25=2225"),
    @("79×99=7821", "47×74=3478"),
    @("76×59=4484", "96×95=9120"),
    @("26×50=1300", "40×82=3280"),
    @("65×48=3120", "72×91=6552"),
    @("80×84=6720", "59×93=5487")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
